# Convert the manual title block into a pandoc-style title block.
#
# Before:
#   P1: "The Eleventh Virgin" (italic) + " " + <line break> + " " +
#       "Part III, Chapter IV ======================="
#   P2: "By Dorothy Day" (bold)
#
# After (single paragraph, no direct formatting):
#   P1: "% Dorothy Day"

$d = $word.ActiveDocument

# Remove the whole first paragraph (title/chapter banner), mark and all, so
# the old "By Dorothy Day" paragraph becomes the new first paragraph.
$firstPara = $d.Paragraphs.Item(1)
$d.Range(0, $firstPara.Range.End).Delete()

# This paragraph used to read "By Dorothy Day" in bold; turn it into a plain
# (non-bold) pandoc title-block author line reading "% Dorothy Day".
$byLine = $d.Paragraphs.Item(1)
$byLineText = $d.Range($byLine.Range.Start, $byLine.Range.End - 1)
$byLineText.Font.Bold = 0

$byLine2 = $d.Paragraphs.Item(1)
$byLineText2 = $d.Range($byLine2.Range.Start, $byLine2.Range.End - 1)
$byLineText2.Text = "% Dorothy Day"
